$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H: "Save" header, formatted like the other header cells (copy G1's style)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the "Save" column values for rows 2-9 (plain numbers, no special style)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 0
